$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = "n.a."
$ws.Range("G7:I7").Value = "n.a."
$ws.Range("E11").Value = "n.a."
$ws.Range("G11:I11").Value = "n.a."
